# Apply scraped-data refresh: increment visit/view counters in column F
# across the four sheets, matching the upstream "gh-pages output" commit.

$wb = $excel.ActiveWorkbook

# 展览 (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 6697
$ws1.Range("F12").Value = 1323
$ws1.Range("F16").Value = 541
$ws1.Range("F21").Value = 716
$ws1.Range("F22").Value = 359
$ws1.Range("F29").Value = 2344
$ws1.Range("F34").Value = 3822

# 演出 (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F6").Value = 744

# 本地生活 (Local Life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F4").Value = 1245
$ws3.Range("F8").Value = 941

# 全部类型 (All Types - aggregate sheet)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 1245
$ws4.Range("F6").Value  = 941
$ws4.Range("F9").Value  = 6697
$ws4.Range("F13").Value = 744
$ws4.Range("F22").Value = 1323
$ws4.Range("F25").Value = 541
$ws4.Range("F29").Value = 716
$ws4.Range("F38").Value = 2344
$ws4.Range("F48").Value = 3822
